$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# hunk 0: @@ -2004,22 +2004,22 @@
$ws_ALC.Range("H28").Value = 203.57143
$ws_ALC.Range("I28").Value = 197.25
$ws_ALC.Range("K28").Value = 197.25
$ws_ALC.Range("M28").Value = 287.75

# hunk 1: @@ -4025,25 +4025,25 @@
$ws_ALC.Range("H69").Value = 3833.3333
$ws_ALC.Range("I69").Value = 2500
$ws_ALC.Range("J69").Value = 4500
$ws_ALC.Range("K69").Value = 7500
$ws_ALC.Range("L69").Value = 13500
$ws_ALC.Range("M69").Value = -6626
$ws_ALC.Range("N69").Value = -15248

# hunk 2: @@ -4175,25 +4175,25 @@
$ws_ALC.Range("H72").Value = 3833.3333
$ws_ALC.Range("I72").Value = 2500
$ws_ALC.Range("J72").Value = 4500
$ws_ALC.Range("K72").Value = 22500
$ws_ALC.Range("L72").Value = 40500
$ws_ALC.Range("M72").Value = -18132
$ws_ALC.Range("N72").Value = -49236

# hunk 3: @@ -4573,25 +4573,25 @@
$ws_ALC.Range("H80").Value = 1704.5454
$ws_ALC.Range("I80").Value = 1750
$ws_ALC.Range("J80").Value = 1694.4445
$ws_ALC.Range("K80").Value = 5250
$ws_ALC.Range("L80").Value = 5083.333500000001
$ws_ALC.Range("M80").Value = -4252
$ws_ALC.Range("N80").Value = -7079.333500000001

# hunk 4: @@ -4720,25 +4720,25 @@
$ws_ALC.Range("H83").Value = 1704.5454
$ws_ALC.Range("I83").Value = 1750
$ws_ALC.Range("J83").Value = 1694.4445
$ws_ALC.Range("K83").Value = 15750
$ws_ALC.Range("L83").Value = 15250.0005
$ws_ALC.Range("M83").Value = -10758
$ws_ALC.Range("N83").Value = -25234.0005

# hunk 5: @@ -4867,25 +4867,25 @@
$ws_ALC.Range("H86").Value = 50970.25
$ws_ALC.Range("I86").Value = 174999
$ws_ALC.Range("J86").Value = 9627.333000000001
$ws_ALC.Range("K86").Value = 174999
$ws_ALC.Range("L86").Value = 9627.333000000001
$ws_ALC.Range("M86").Value = -173876
$ws_ALC.Range("N86").Value = -11873.333

# hunk 6: @@ -4919,22 +4919,22 @@
$ws_ALC.Range("H87").Value = 39993.332
$ws_ALC.Range("J87").Value = 39993.332
$ws_ALC.Range("L87").Value = 39993.332
$ws_ALC.Range("N87").Value = -42489.332

# hunk 7: @@ -5020,25 +5020,25 @@
$ws_ALC.Range("H89").Value = 50970.25
$ws_ALC.Range("I89").Value = 174999
$ws_ALC.Range("J89").Value = 9627.333000000001
$ws_ALC.Range("K89").Value = 874995
$ws_ALC.Range("L89").Value = 48136.665
$ws_ALC.Range("M89").Value = -869379
$ws_ALC.Range("N89").Value = -59368.665

# hunk 8: @@ -5072,22 +5072,22 @@
$ws_ALC.Range("H90").Value = 39993.332
$ws_ALC.Range("J90").Value = 39993.332
$ws_ALC.Range("L90").Value = 119979.996
$ws_ALC.Range("N90").Value = -132459.996

# hunk 9: @@ -5871,22 +5871,22 @@
$ws_ALC.Range("H106").Value = 3395.4443
$ws_ALC.Range("I106").Value = 3395.4443
$ws_ALC.Range("K106").Value = 3395.4443
$ws_ALC.Range("M106").Value = -2764.4443

# hunk 10: @@ -6162,25 +6162,25 @@
$ws_ALC.Range("H112").Value = 2407.6924
$ws_ALC.Range("J112").Value = 2770
$ws_ALC.Range("L112").Value = 8310
$ws_ALC.Range("N112").Value = -10526

# hunk 11: @@ -7574,22 +7574,22 @@
$ws_ALC.Range("H141").Value = 2293.1555
$ws_ALC.Range("I141").Value = 1204.8422
$ws_ALC.Range("K141").Value = 3614.5266
$ws_ALC.Range("M141").Value = 1565.4734

# hunk 12: @@ -7772,22 +7772,22 @@
$ws_ARM.Range("H3").Value = 505
$ws_ARM.Range("I3").Value = 505
$ws_ARM.Range("J3").Value = 0
$ws_ARM.Range("K3").Value = 505
$ws_ARM.Range("L3").ClearContents()   # removed cell
$ws_ARM.Range("M3").Value = -390   # new cell
$ws_ARM.Range("N3").Value = 0

# hunk 13: @@ -9842,22 +9842,22 @@
$ws_ARM.Range("H45").Value = 2328.6
$ws_ARM.Range("I45").Value = 2269.25
$ws_ARM.Range("K45").Value = 2269.25
$ws_ARM.Range("M45").Value = -1892.25

# hunk 14: @@ -10280,22 +10280,22 @@
$ws_ARM.Range("H54").Value = 32000
$ws_ARM.Range("J54").Value = 32000
$ws_ARM.Range("L54").Value = 32000
$ws_ARM.Range("N54").Value = -33538

# hunk 15: @@ -12994,22 +12994,22 @@
$ws_ARM.Range("H110").Value = 409.6
$ws_ARM.Range("I110").Value = 409.6
$ws_ARM.Range("K110").Value = 409.6
$ws_ARM.Range("M110").Value = 1635.4

# hunk 16: @@ -16408,19 +16408,22 @@
$ws_BSM.Range("H38").Value = 30000
$ws_BSM.Range("J38").Value = 30000
$ws_BSM.Range("L38").Value = 30000
$ws_BSM.Range("N38").Value = -30832   # new cell

# hunk 17: @@ -19095,25 +19098,25 @@
$ws_BSM.Range("H94").Value = 1769.85
$ws_BSM.Range("I94").Value = 1400.5
$ws_BSM.Range("J94").Value = 2323.875
$ws_BSM.Range("K94").Value = 1400.5
$ws_BSM.Range("L94").Value = 2323.875
$ws_BSM.Range("M94").Value = -949.5
$ws_BSM.Range("N94").Value = -3225.875

# hunk 18: @@ -19340,22 +19343,22 @@
$ws_BSM.Range("H99").Value = 2499.6667
$ws_BSM.Range("I99").Value = 2272.9092
$ws_BSM.Range("K99").Value = 2272.9092
$ws_BSM.Range("M99").Value = -774.9092000000001

# hunk 19: @@ -22099,22 +22102,22 @@
$ws_CRP.Range("H14").Value = 12407.6
$ws_CRP.Range("J14").Value = 12407.6
$ws_CRP.Range("L14").Value = 12407.6
$ws_CRP.Range("N14").Value = -12747.6

# hunk 20: @@ -22941,25 +22944,25 @@
$ws_CRP.Range("H31").Value = 1220.8667
$ws_CRP.Range("I31").Value = 972.8182
$ws_CRP.Range("J31").Value = 1903
$ws_CRP.Range("K31").Value = 972.8182
$ws_CRP.Range("L31").Value = 1903
$ws_CRP.Range("M31").Value = -677.8182
$ws_CRP.Range("N31").Value = -2493

# hunk 21: @@ -23097,25 +23100,25 @@
$ws_CRP.Range("H34").Value = 1220.8667
$ws_CRP.Range("I34").Value = 972.8182
$ws_CRP.Range("J34").Value = 1903
$ws_CRP.Range("K34").Value = 972.8182
$ws_CRP.Range("L34").Value = 1903
$ws_CRP.Range("M34").Value = -770.8182
$ws_CRP.Range("N34").Value = -2307

# hunk 22: @@ -25603,22 +25606,22 @@
$ws_CRP.Range("H86").Value = 10099.6
$ws_CRP.Range("I86").Value = 7000
$ws_CRP.Range("K86").Value = 7000
$ws_CRP.Range("M86").Value = -5877

# hunk 23: @@ -25750,22 +25753,22 @@
$ws_CRP.Range("H89").Value = 10099.6
$ws_CRP.Range("I89").Value = 7000
$ws_CRP.Range("K89").Value = 35000
$ws_CRP.Range("M89").Value = -29384

# hunk 24: @@ -26096,22 +26099,22 @@
$ws_CRP.Range("H96").Value = 16312
$ws_CRP.Range("J96").Value = 16312
$ws_CRP.Range("L96").Value = 16312
$ws_CRP.Range("N96").Value = -21804

# hunk 25: @@ -26620,22 +26623,22 @@
$ws_CRP.Range("H107").Value = 535.1667
$ws_CRP.Range("I107").Value = 460.07144
$ws_CRP.Range("K107").Value = 460.07144
$ws_CRP.Range("M107").Value = 1459.92856

# hunk 26: @@ -27812,22 +27815,22 @@
$ws_CRP.Range("H132").Value = 3642.5
$ws_CRP.Range("I132").Value = 3642.5
$ws_CRP.Range("K132").Value = 10927.5
$ws_CRP.Range("M132").Value = -8397.5

# hunk 27: @@ -27910,25 +27913,25 @@
$ws_CRP.Range("H134").Value = 4575
$ws_CRP.Range("J134").Value = 4496
$ws_CRP.Range("L134").Value = 13488
$ws_CRP.Range("N134").Value = -18558

# hunk 28: @@ -28497,22 +28500,22 @@
$ws_CUL.Range("H4").Value = 1301.85
$ws_CUL.Range("I4").Value = 1255.1177
$ws_CUL.Range("K4").Value = 3765.3531
$ws_CUL.Range("M4").Value = -3653.3531

# hunk 29: @@ -30248,22 +30251,22 @@
$ws_CUL.Range("H39").Value = 2053.9
$ws_CUL.Range("J39").Value = 2053.9
$ws_CUL.Range("L39").Value = 6161.700000000001
$ws_CUL.Range("N39").Value = -6749.700000000001

# hunk 30: @@ -40388,25 +40391,25 @@
$ws_GSM.Range("H105").Value = 32017.4
$ws_GSM.Range("I105").Value = 8000
$ws_GSM.Range("J105").Value = 53032.625
$ws_GSM.Range("K105").Value = 8000
$ws_GSM.Range("L105").Value = 53032.625
$ws_GSM.Range("M105").Value = -4506
$ws_GSM.Range("N105").Value = -60020.625

# hunk 31: @@ -44094,22 +44097,22 @@
$ws_LTW.Range("H40").Value = 8316
$ws_LTW.Range("I40").Value = 7449
$ws_LTW.Range("K40").Value = 7449
$ws_LTW.Range("M40").Value = -7313

# hunk 32: @@ -44802,25 +44805,25 @@
$ws_LTW.Range("H55").Value = 308.45456
$ws_LTW.Range("I55").Value = 293.2
$ws_LTW.Range("J55").Value = 321.16666
$ws_LTW.Range("K55").Value = 293.2
$ws_LTW.Range("L55").Value = 321.16666
$ws_LTW.Range("M55").Value = -120.2
$ws_LTW.Range("N55").Value = -667.16666

# hunk 33: @@ -48660,22 +48663,22 @@
$ws_LTW.Range("H136").Value = 3516.8333
$ws_LTW.Range("I136").Value = 3211
$ws_LTW.Range("K136").Value = 9633
$ws_LTW.Range("M136").Value = -7083

# hunk 34: @@ -51970,25 +51973,25 @@
$ws_WVR.Range("H62").Value = 3603.6667
$ws_WVR.Range("J62").Value = 1666.6666
$ws_WVR.Range("L62").Value = 1666.6666
$ws_WVR.Range("N62").Value = -2914.6666

# hunk 35: @@ -52114,25 +52117,25 @@
$ws_WVR.Range("H65").Value = 3603.6667
$ws_WVR.Range("J65").Value = 1666.6666
$ws_WVR.Range("L65").Value = 8333.333000000001
$ws_WVR.Range("N65").Value = -14573.333

# hunk 36: @@ -55545,25 +55548,25 @@
$ws_WVR.Range("H136").Value = 7202
$ws_WVR.Range("I136").Value = 7948.5
$ws_WVR.Range("J136").Value = 6642.125
$ws_WVR.Range("K136").Value = 23845.5
$ws_WVR.Range("L136").Value = 19926.375
$ws_WVR.Range("M136").Value = -21295.5
# N136: unchanged (-25026.375) - skip
